$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Prepay Loan" sheet: correct the prepay amount (10145.16 -> 10045.16) and
# move the selection to B14 (also drops this sheet's tabSelected flag once
# another sheet is activated below).
# ---------------------------------------------------------------------------
$wsPrepay = $wb.Worksheets.Item("Prepay Loan")
$wsPrepay.Range("B4").Value = 10045.16
$wsPrepay.Activate()
$wsPrepay.Range("B14").Select()

# ---------------------------------------------------------------------------
# "Repayment schedule" sheet: column widths were narrowed (columns A, C, I,
# and K:L got tightened - K and L now share a single uniform width) and the
# selection moved to M10.
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns.Item(1).ColumnWidth = 1.1666666666666667
$wsSchedule.Columns.Item(3).ColumnWidth = 8.5
$wsSchedule.Columns.Item(9).ColumnWidth = 4.3333333333333333
$wsSchedule.Range($wsSchedule.Columns.Item(11), $wsSchedule.Columns.Item(12)).ColumnWidth = 8.3333333333333333
$wsSchedule.Activate()
$wsSchedule.Range("M10").Select()

# ---------------------------------------------------------------------------
# "Transactions" sheet: the three oldest "Accrual" transactions (rows 2-4)
# are removed; the remaining Repayment (was row 5) and Disbursement (was row
# 6) rows shift up to become rows 2 and 3. The Repayment row's total amount
# is also corrected from 10145.16 to 10045.16. This sheet becomes the active
# tab/selected sheet, with the cursor on G7.
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Rows("2:4").Delete()
$wsTransactions.Range("E2").Value = 10045.16
$wsTransactions.Activate()
$wsTransactions.Range("G7").Select()
